$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.Formula = ('="' + $text + '"')
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$excel.CutCopyMode = $false

# --- Step 1: prepare the red-fill style (same as G6/G16/G22) on I5, I7, I22 before stripping it elsewhere ---
$ws.Cells.Item(6,7).Copy() | Out-Null
$ws.Cells.Item(5,9).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(6,7).Copy() | Out-Null
$ws.Cells.Item(7,9).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(6,7).Copy() | Out-Null
$ws.Cells.Item(22,9).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Step 2: header cells H1, I1 get the same style as G1 (bold/border/center) ---
$ws.Cells.Item(1,7).Copy() | Out-Null
$ws.Cells.Item(1,8).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1,7).Copy() | Out-Null
$ws.Cells.Item(1,9).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
Set-TextValue $ws.Cells.Item(1,8) "2025-09-16"
Set-TextValue $ws.Cells.Item(1,9) "2025-09-18"

# --- Step 3: fill in attendance values for columns H (2025-09-16) and I (2025-09-18) ---
Set-TextValue $ws.Cells.Item(2,8) "0"
Set-TextValue $ws.Cells.Item(2,9) "0"
Set-TextValue $ws.Cells.Item(3,8) "0"
Set-TextValue $ws.Cells.Item(3,9) "0"
Set-TextValue $ws.Cells.Item(4,8) "0"
Set-TextValue $ws.Cells.Item(4,9) "0"
Set-TextValue $ws.Cells.Item(5,8) "1"
Set-TextValue $ws.Cells.Item(5,9) "1"
Set-TextValue $ws.Cells.Item(6,8) "0"
Set-TextValue $ws.Cells.Item(6,9) "0"
Set-TextValue $ws.Cells.Item(7,8) "0"
Set-TextValue $ws.Cells.Item(7,9) "1"
Set-TextValue $ws.Cells.Item(8,8) "0"
Set-TextValue $ws.Cells.Item(8,9) "0"
Set-TextValue $ws.Cells.Item(9,8) "0"
Set-TextValue $ws.Cells.Item(9,9) "0"
Set-TextValue $ws.Cells.Item(10,8) "0"
Set-TextValue $ws.Cells.Item(10,9) "0"
Set-TextValue $ws.Cells.Item(11,8) "0"
Set-TextValue $ws.Cells.Item(11,9) "0"
Set-TextValue $ws.Cells.Item(12,8) "0"
Set-TextValue $ws.Cells.Item(12,9) "0"
Set-TextValue $ws.Cells.Item(13,8) "0"
Set-TextValue $ws.Cells.Item(13,9) "0"
Set-TextValue $ws.Cells.Item(14,8) "0"
Set-TextValue $ws.Cells.Item(14,9) "0"
Set-TextValue $ws.Cells.Item(15,8) "0.5"
Set-TextValue $ws.Cells.Item(15,9) "0"
Set-TextValue $ws.Cells.Item(16,8) "0"
Set-TextValue $ws.Cells.Item(16,9) "0"
Set-TextValue $ws.Cells.Item(17,8) "0"
Set-TextValue $ws.Cells.Item(17,9) "0"
Set-TextValue $ws.Cells.Item(18,8) "0"
Set-TextValue $ws.Cells.Item(18,9) "0"
Set-TextValue $ws.Cells.Item(19,8) "0"
Set-TextValue $ws.Cells.Item(19,9) "0"
Set-TextValue $ws.Cells.Item(20,8) "0"
Set-TextValue $ws.Cells.Item(20,9) "0"
Set-TextValue $ws.Cells.Item(21,8) "0"
Set-TextValue $ws.Cells.Item(21,9) "0"
Set-TextValue $ws.Cells.Item(22,8) "0"
$ws.Cells.Item(22,9).Value = 1
Set-TextValue $ws.Cells.Item(23,8) "0"
Set-TextValue $ws.Cells.Item(23,9) "0"
Set-TextValue $ws.Cells.Item(24,8) "0"
Set-TextValue $ws.Cells.Item(24,9) "0"
Set-TextValue $ws.Cells.Item(25,8) "0"
Set-TextValue $ws.Cells.Item(25,9) "0"
Set-TextValue $ws.Cells.Item(26,8) "0"
Set-TextValue $ws.Cells.Item(26,9) "0"
Set-TextValue $ws.Cells.Item(27,8) "0"
Set-TextValue $ws.Cells.Item(27,9) "0"
Set-TextValue $ws.Cells.Item(28,8) "0"
Set-TextValue $ws.Cells.Item(28,9) "0"
Set-TextValue $ws.Cells.Item(29,8) "0"
Set-TextValue $ws.Cells.Item(29,9) "0"
Set-TextValue $ws.Cells.Item(30,8) "0"
Set-TextValue $ws.Cells.Item(30,9) "0"
Set-TextValue $ws.Cells.Item(31,8) "0"
Set-TextValue $ws.Cells.Item(31,9) "0"

# --- Step 4: remove the now-stale highlight styles from column G ---
$ws.Cells.Item(6,7).Style = "Normal"
$ws.Cells.Item(8,7).Style = "Normal"
$ws.Cells.Item(16,7).Style = "Normal"
$ws.Cells.Item(22,7).Style = "Normal"
Set-TextValue $ws.Cells.Item(27,7) "1"
$ws.Cells.Item(27,7).Style = "Normal"
